# "Tests with width search" - fill in the remaining widthSearch_v2 runs
# (ITER_LIMIT 1000000; ITER_LIMIT 100000/MEMORY_SIZE 500; ITER_LIMIT 100000/MEMORY_SIZE 1000)
# on the "Largura" sheet, and move the selection on both result sheets.

$wb = $excel.ActiveWorkbook

$largura = $wb.Worksheets.Item("Largura")

# Row 11: ITER_LIMIT=1000000
$largura.Range("E11").Value = "N"
$largura.Range("G11").Value = 859.178
$largura.Range("H11").Value = 1813431
$largura.Range("I11").Value = 17

# Row 13: ITER_LIMIT=100000, MEMORY_SIZE=500
$largura.Range("E13").Value = "N"
$largura.Range("G13").Value = 246691
$largura.Range("H13").Value = 187086
$largura.Range("I13").Value = 17

# Row 14: ITER_LIMIT=100000, MEMORY_SIZE=1000
$largura.Range("E14").Value = "N"
$largura.Range("G14").Value = 498.404
$largura.Range("G14").NumberFormat = "#,##0.000"
$largura.Range("H14").Value = 181768
$largura.Range("I14").Value = 18

# Update the saved cursor/selection on both sheets touched in this session.
$profundidade = $wb.Worksheets.Item("Profundidade")
$profundidade.Range("N1").Select() | Out-Null

$largura.Range("E12").Select() | Out-Null
